$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") holds a comma-separated list of actors who
# touched each attendance record. Re-order each list so that any
# "System"/"system" entries come first (stable order among themselves),
# followed by the remaining entries in their original relative order.
# When there is no System-ish entry at all, the list is simply reversed.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ',\s*'
    if ($parts.Count -le 1) { continue }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq 'system') {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -gt 0) {
        $newParts = $systemParts + $otherParts
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newText = [string]::Join(', ', $newParts)
    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
